$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Filtro")
[void]$ws.Activate()
$ws.Range("B8").Value = 647
[void]$ws.Range("D18").Select()
